# Update crypto price (D) and 1h volume change (E) columns with latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.268.90"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.501.07"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "4.098.74"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "3.501.74"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.49%  "
$ws.Range("D17").Value = "64.368.71"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.571"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "3.642.33"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.02%  "
$ws.Range("D33").Value = "3.524.55"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.60%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("E47").Value = "  -4.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").Value = "2.460.89"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.892"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.95%  "
